$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the query table (ListObject) with 5 additional rows so the
# table / sheet range grows from A1:F28 to A1:F33, matching the refreshed
# SharePoint list export.
$lo = $ws.ListObjects.Item(1)
for ($i = 0; $i -lt 5; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# Ensure the new rows use the same Text number format as the rest of
# column A / E / F (so formatting matches the other data rows).
$ws.Range("A29:A33").NumberFormat = "@"
$ws.Range("E29:E33").NumberFormat = "@"
$ws.Range("F29:F33").NumberFormat = "@"

# Re-write the Software / Related Business Need / Related Deliverable
# columns for every data row so the table reflects the refreshed list
# (renamed items, corrected labels and newly added software entries).
$ws.Range("A2").Value = '_Starting Over'
$ws.Range("B2").Value = '_Starting Over;#32;#Brainstorming;#33;#Collaboration;#34;#Communication;#42;#Coordination;#43;#Decision Making;#36;#Evaluation;#44;#Exchange;#45;#Feedback;#35;#Learning;#37;#Referencing;#38;#Research;#39;#Simulation;#40;#Sysnthesis;#46;#Videography;#41'
$ws.Range("C2").Value = '_Starting Over;#22;#Activity Coordination;#53;#Blog;#30;#Change Management;#54;#Citation and Bibliographies;#40;#Course Coordination;#55;#Curriculum;#34;#Discussion Forum;#29;#Document Collaboration;#25;#eBook;#35;#Evaluation Plan;#56;#HTA;#58;#Ideation;#26;#InfoGraphic;#47;#Job Aid;#36;#Knowledge Map;#24;#Knowledge Repositories;#43;#Lessons Learned Repository;#31;#Literature Review;#44;#Memes and Vemes;#50;#Newsletter;#49;#Online Module;#37;#Posters;#48;#Presentations;#38;#Print Screens;#46;#Project Coordination;#52;#Qualitative and Quantitative Analysis;#42;#Reporting Tool;#33;#Requests Coordination;#51;#ROI;#57;#Subject Matter Expertise;#27;#Survey;#28;#Toolbox: Methods, Methodologies, Tools;#32;#Training Manual;#39;#Video;#45'
$ws.Range("A3").Value = 'Adobe'
$ws.Range("B3").Value = '_Starting Over;#32;#Simulation;#40'
$ws.Range("C3").Value = '_Starting Over;#22;#Video;#45'
$ws.Range("A4").Value = 'Adobe Acrobat Pro'
$ws.Range("B4").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C4").Value = '_Starting Over;#22;#eBook;#35;#Information Sheet;#60'
$ws.Range("A5").Value = 'Adobe Audition'
$ws.Range("B5").Value = '_Starting Over;#32;#Audio Recording;#47'
$ws.Range("C5").Value = '_Starting Over;#22;#Audio;#61'
$ws.Range("A6").Value = 'Adobe eLearning'
$ws.Range("B6").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C6").Value = '_Starting Over;#22;#Online Module;#37;#Training Manual;#39'
$ws.Range("A7").Value = 'Adobe Lifecycle'
$ws.Range("B7").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C7").Value = '_Starting Over;#22;#Information Sheet;#60'
$ws.Range("A8").Value = 'Camtasia'
$ws.Range("B8").Value = '_Starting Over;#32;#Simulation;#40'
$ws.Range("C8").Value = '_Starting Over;#22;#Video;#45'
$ws.Range("A9").Value = 'Cmap'
$ws.Range("B9").Value = 'Brainstorming;#33;#_Starting Over;#32'
$ws.Range("C9").Value = 'Knowledge Map;#24;#_Starting Over;#22'
$ws.Range("A10").Value = 'Comapping'
$ws.Range("B10").Value = 'Brainstorming;#33;#_Starting Over;#32'
$ws.Range("C10").Value = 'Knowledge Map;#24;#_Starting Over;#22'
$ws.Range("A11").Value = 'Designer ES2'
$ws.Range("B11").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C11").Value = '_Starting Over;#22;#Information Sheet;#60'
$ws.Range("A12").Value = 'Go Animate'
$ws.Range("B12").Value = '_Starting Over;#32;#Videography;#41'
$ws.Range("C12").Value = '_Starting Over;#22;#Video;#45'
$ws.Range("A13").Value = 'KRS Lib Guides'
$ws.Range("B13").Value = '_Starting Over;#32;#Research;#39'
$ws.Range("C13").Value = '_Starting Over;#22;#Literature Review;#44'
$ws.Range("A14").Value = 'Lync'
$ws.Range("B14").Value = 'Collaboration;#34;#_Starting Over;#32'
$ws.Range("C14").Value = 'Document Collaboration;#25;#_Starting Over;#22'
$ws.Range("A15").Value = 'Lync Annotation Board'
$ws.Range("B15").Value = '_Starting Over;#32;#Brainstorming;#33'
$ws.Range("C15").Value = '_Starting Over;#22;#Knowledge Map;#24'
$ws.Range("A16").Value = 'Mendeley'
$ws.Range("B16").Value = '_Starting Over;#32;#Referencing;#38'
$ws.Range("C16").Value = '_Starting Over;#22;#Citation and Bibliographies;#40'
$ws.Range("A17").Value = 'Microsoft Access'
$ws.Range("B17").Value = '_Starting Over;#32;#Decision Making;#36'
$ws.Range("C17").Value = '_Starting Over;#22;#Reporting Tool;#33'
$ws.Range("A18").Value = 'Microsoft PowerPoint'
$ws.Range("B18").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C18").Value = '_Starting Over;#22;#Presentations;#38'
$ws.Range("A19").Value = 'Microsoft Publisher'
$ws.Range("B19").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C19").Value = '_Starting Over;#22;#Information Sheet;#60'
$ws.Range("A20").Value = 'Microsoft Word'
$ws.Range("B20").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C20").Value = '_Starting Over;#22;#Curriculum;#34;#Job Aid;#36'
$ws.Range("A21").Value = 'Mindomo Mind Map'
$ws.Range("B21").Value = 'Brainstorming;#33;#_Starting Over;#32'
$ws.Range("C21").Value = 'Knowledge Map;#24;#_Starting Over;#22'
$ws.Range("A22").Value = 'nVivo'
$ws.Range("B22").Value = '_Starting Over;#32;#Research;#39'
$ws.Range("C22").Value = '_Starting Over;#22;#Qualitative and Quantitative Analysis;#42'
$ws.Range("A23").Value = 'Optimal Workshop'
$ws.Range("B23").Value = 'Brainstorming;#33;#_Starting Over;#32'
$ws.Range("C23").Value = 'Knowledge Map;#24;#_Starting Over;#22'
$ws.Range("A24").Value = 'Presenter Media'
$ws.Range("B24").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C24").Value = '_Starting Over;#22;#Presentations;#38'
$ws.Range("A25").Value = 'Prespectore'
$ws.Range("B25").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C25").Value = '_Starting Over;#22;#Presentations;#38'
$ws.Range("A26").Value = 'Prezi'
$ws.Range("B26").Value = '_Starting Over;#32;#Learning;#37'
$ws.Range("C26").Value = '_Starting Over;#22;#Presentations;#38'
$ws.Range("A27").Value = 'RefWorks'
$ws.Range("B27").Value = '_Starting Over;#32;#Referencing;#38'
$ws.Range("C27").Value = '_Starting Over;#22;#Citation and Bibliographies;#40'
$ws.Range("A28").Value = 'SharePoint'
$ws.Range("B28").Value = 'Collaboration;#34;#Feedback;#35;#_Starting Over;#32;#Research;#39;#Coordination;#43'
$ws.Range("C28").Value = 'Document Collaboration;#25;#Subject Matter Expertise;#27;#Discussion Forum;#29;#Blog;#30;#Lessons Learned Repository;#31;#_Starting Over;#22;#Toolbox: Methods, Methodologies, Tools;#32;#Reporting Tool;#33;#Knowledge Repositories;#43;#Requests Coordination;#51;#Project Coordination;#52;#Activity Coordination;#53;#Change Management;#54;#Environmental Scan;#59'
$ws.Range("A29").Value = 'Snagit'
$ws.Range("B29").Value = '_Starting Over;#32;#Simulation;#40'
$ws.Range("C29").Value = '_Starting Over;#22;#Video;#45;#Print Screens;#46'
$ws.Range("A30").Value = 'Soney Vegas'
$ws.Range("B30").Value = '_Starting Over;#32;#Videography;#41'
$ws.Range("C30").Value = '_Starting Over;#22;#Video;#45'
$ws.Range("A31").Value = 'Statit'
$ws.Range("B31").Value = '_Starting Over;#32;#Decision Making;#36'
$ws.Range("C31").Value = '_Starting Over;#22;#Reporting Tool;#33'
$ws.Range("A32").Value = 'Survey Net'
$ws.Range("B32").Value = 'Feedback;#35;#_Starting Over;#32'
$ws.Range("C32").Value = 'Survey;#28;#_Starting Over;#22'
$ws.Range("A33").Value = 'Tableau'
$ws.Range("B33").Value = '_Starting Over;#32;#Decision Making;#36'
$ws.Range("C33").Value = '_Starting Over;#22;#Reporting Tool;#33'

# Fill in Item Type / Path for the newly-added rows (existing rows
# already carry these values).
for ($r = 29; $r -le 33; $r++) {
    $ws.Range("E" + $r).Value = "Item"
    $ws.Range("F" + $r).Value = "teams/kmqa/Lists/Software"
}

# Update the hidden defined name that the query table range uses so it
# also reflects the new F33 extent.
$wb.Names.Item(1).RefersTo = "=owssvr!`$A`$1:`$F`$33"

Write-Host "done"
